$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (order submission) ---
$ws.Range("A2").Value = "ΑΠΟΣΤΟΛΟΣ ΑΝΑΣΤΑΣΙΟΥ"

# B2/C2 hold numeric-looking codes that must stay text, force text format first
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "6362"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "6363"

# D2 (phone #1) stays the same: 6975362321

# Clear out phone #2, postal code, area and address
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""

$ws.Range("L2").Value = "DTS-0019"
$ws.Range("M2").Value = "ΣΥΜΒΑΤΗ ΜΕΛΑΝΟΤΑΙΝΙΑ EPSON ERC-32B BLACK"
$ws.Range("N2").Value = 15

# O2 ("20€") must remain text so the currency sign survives
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "20€"

# --- Update row 3 (second order line) ---
$ws.Range("L3").Value = "GPI-0137"
$ws.Range("M3").Value = "BROTHER INK LC1240 CYAN ΣΥΜΒΑΤΟ 10ml"
$ws.Range("N3").Value = 2
# O3 stays empty

# --- Remove rows 4-6 entirely (only two order lines remain) ---
$ws.Range("A4:O6").EntireRow.Delete()
